$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells B1:G1 from A..F to rating1..rating6
$ws.Range("B1").Value = "rating1"
$ws.Range("C1").Value = "rating2"
$ws.Range("D1").Value = "rating3"
$ws.Range("E1").Value = "rating4"
$ws.Range("F1").Value = "rating5"
$ws.Range("G1").Value = "rating6"

# Clear the stray "good"/"fine"/"x"/"."/"10a"/"10b"/"2/4a"/"2/4b" values
$ws.Range("B3:G3").ClearContents()
$ws.Range("B8:G8").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A102").ClearContents()
$ws.Range("A103").ClearContents()
$ws.Range("A104").ClearContents()

# Select B1:G1 (matches the saved cursor position)
[void]$ws.Range("B1:G1").Select()

# Turn on AutoFilter for the data range and register the hidden
# sheet-scoped _FilterDatabase defined name that Excel writes alongside it.
[void]$ws.Range("A1:G125").AutoFilter()
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$125")
$name.Visible = $false
